$wb = $excel.ActiveWorkbook

# --- Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Cells.Item(4, 2).Value = "inf"
$ws.Cells.Item(6, 2).Value = 1161943.262425547
$ws.Cells.Item(8, 2).Value = 492028.9342484446
$ws.Cells.Item(10, 2).Value = 4536306.814155112

# --- Costs and Revenues ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Cells.Item(2, 5).Value = 521718.6705202902
$ws.Cells.Item(2, 7).Value = 523898.4668372231
$ws.Cells.Item(2, 8).Value = 523898.4668372231
$ws.Cells.Item(2, 9).Value = 523898.4668372231
$ws.Cells.Item(2, 10).Value = 521718.6705202902
$ws.Cells.Item(2, 12).Value = 521718.6705202902
$ws.Cells.Item(2, 13).Value = 523898.4668372231
$ws.Cells.Item(2, 14).Value = 523898.4668372231
$ws.Cells.Item(2, 15).Value = 521718.6705202902
$ws.Cells.Item(2, 16).Value = 521718.6705202902
$ws.Cells.Item(3, 5).Value = 317261.5624159134
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 3383.95761854822
$ws.Cells.Item(3, 10).Value = 314523.1248318268
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 6122.395202634834
$ws.Cells.Item(3, 15).Value = 311784.6872477402
$ws.Cells.Item(4, 5).Value = 297779.8398421858
$ws.Cells.Item(4, 7).Value = 299023.995811336
$ws.Cells.Item(4, 8).Value = 299023.995811336
$ws.Cells.Item(4, 9).Value = 299023.995811336
$ws.Cells.Item(4, 10).Value = 297779.8398421858
$ws.Cells.Item(4, 12).Value = 297779.8398421858
$ws.Cells.Item(4, 13).Value = 299023.995811336
$ws.Cells.Item(4, 14).Value = 299023.995811336
$ws.Cells.Item(4, 15).Value = 297779.8398421858
$ws.Cells.Item(4, 16).Value = 297779.8398421858
$ws.Cells.Item(5, 5).Value = 33407.66314315072
$ws.Cells.Item(5, 7).Value = 33695.43527972145
$ws.Cells.Item(5, 8).Value = 33695.43527972145
$ws.Cells.Item(5, 9).Value = 33695.43527972145
$ws.Cells.Item(5, 10).Value = 33407.66314315072
$ws.Cells.Item(5, 12).Value = 33407.66314315072
$ws.Cells.Item(5, 13).Value = 33695.43527972145
$ws.Cells.Item(5, 14).Value = 33695.43527972145
$ws.Cells.Item(5, 15).Value = 33407.66314315072
$ws.Cells.Item(5, 16).Value = 33407.66314315072
$ws.Cells.Item(6, 5).Value = -126730.3948809597
$ws.Cells.Item(6, 6).Value = 190531.1675349537
$ws.Cells.Item(6, 7).Value = 187795.0781276174
$ws.Cells.Item(6, 8).Value = 191179.0357461657
$ws.Cells.Item(6, 9).Value = 191179.0357461657
$ws.Cells.Item(6, 10).Value = -123991.9572968731
$ws.Cells.Item(6, 11).Value = 190531.1675349537
$ws.Cells.Item(6, 12).Value = 185056.6405435308
$ws.Cells.Item(6, 13).Value = 191179.0357461657
$ws.Cells.Item(6, 14).Value = 191179.0357461657
$ws.Cells.Item(6, 15).Value = -121253.5197127866
$ws.Cells.Item(6, 16).Value = 190531.1675349537

# --- Installed Capacities ---
$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Cells.Item(2, 5).Value = 397.3838530629687
$ws.Cells.Item(2, 7).Value = 400.806900043077
$ws.Cells.Item(2, 8).Value = 400.806900043077
$ws.Cells.Item(2, 9).Value = 400.806900043077
$ws.Cells.Item(2, 10).Value = 397.3838530629687
$ws.Cells.Item(2, 12).Value = 397.3838530629687
$ws.Cells.Item(2, 13).Value = 400.806900043077
$ws.Cells.Item(2, 14).Value = 400.806900043077
$ws.Cells.Item(2, 15).Value = 397.3838530629687
$ws.Cells.Item(2, 16).Value = 397.3838530629687

# --- Added Capacities ---
$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Cells.Item(2, 5).Value = 396.5769530198917
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 4.229947023185275
$ws.Cells.Item(2, 10).Value = 393.1539060397835
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 7.652994003293543
$ws.Cells.Item(2, 15).Value = 389.7308590596753

# --- Retired Capacities ---
$ws = $wb.Worksheets.Item("Retired Capacities")
$ws.Cells.Item(2, 10).Value = 396.5769530198917
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 4.229947023185275
$ws.Cells.Item(2, 15).Value = 393.1539060397835
$ws.Cells.Item(2, 16).Value = 0

# --- DG Dispatch ---
$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Cells.Item(11, 2).Value = 397.3838530629687
$ws.Cells.Item(11, 3).Value = 397.3838530629687
$ws.Cells.Item(11, 4).Value = 397.3838530629687
$ws.Cells.Item(11, 5).Value = 397.3838530629687
$ws.Cells.Item(11, 6).Value = 397.3838530629687
$ws.Cells.Item(11, 7).Value = 397.3838530629687
$ws.Cells.Item(11, 22).Value = 358.9907805655117
$ws.Cells.Item(11, 23).Value = 397.3838530629687
$ws.Cells.Item(11, 24).Value = 397.3838530629687
$ws.Cells.Item(11, 25).Value = 397.3838530629687
$ws.Cells.Item(17, 2).Value = 400.806900043077
$ws.Cells.Item(17, 3).Value = 400.806900043077
$ws.Cells.Item(17, 4).Value = 400.806900043077
$ws.Cells.Item(17, 5).Value = 400.806900043077
$ws.Cells.Item(17, 6).Value = 400.806900043077
$ws.Cells.Item(17, 24).Value = 400.806900043077
$ws.Cells.Item(17, 25).Value = 400.806900043077
$ws.Cells.Item(20, 2).Value = 400.806900043077
$ws.Cells.Item(20, 3).Value = 400.806900043077
$ws.Cells.Item(20, 4).Value = 400.806900043077
$ws.Cells.Item(20, 5).Value = 400.806900043077
$ws.Cells.Item(20, 6).Value = 400.806900043077
$ws.Cells.Item(20, 24).Value = 400.806900043077
$ws.Cells.Item(20, 25).Value = 400.806900043077
$ws.Cells.Item(23, 2).Value = 400.806900043077
$ws.Cells.Item(23, 3).Value = 400.806900043077
$ws.Cells.Item(23, 4).Value = 400.806900043077
$ws.Cells.Item(23, 5).Value = 400.806900043077
$ws.Cells.Item(23, 6).Value = 400.806900043077
$ws.Cells.Item(23, 24).Value = 400.806900043077
$ws.Cells.Item(23, 25).Value = 400.806900043077
$ws.Cells.Item(26, 2).Value = 397.3838530629687
$ws.Cells.Item(26, 3).Value = 397.3838530629687
$ws.Cells.Item(26, 4).Value = 397.3838530629687
$ws.Cells.Item(26, 5).Value = 397.3838530629687
$ws.Cells.Item(26, 6).Value = 397.3838530629687
$ws.Cells.Item(26, 7).Value = 397.3838530629687
$ws.Cells.Item(26, 22).Value = 358.9907805655117
$ws.Cells.Item(26, 23).Value = 397.3838530629687
$ws.Cells.Item(26, 24).Value = 397.3838530629687
$ws.Cells.Item(26, 25).Value = 397.3838530629687
$ws.Cells.Item(29, 2).Value = 397.3838530629687
$ws.Cells.Item(29, 3).Value = 397.3838530629687
$ws.Cells.Item(29, 4).Value = 397.3838530629687
$ws.Cells.Item(29, 5).Value = 397.3838530629687
$ws.Cells.Item(29, 6).Value = 397.3838530629687
$ws.Cells.Item(29, 23).Value = 397.3838530629687
$ws.Cells.Item(29, 24).Value = 397.3838530629687
$ws.Cells.Item(29, 25).Value = 397.3838530629687
$ws.Cells.Item(32, 2).Value = 400.806900043077
$ws.Cells.Item(32, 3).Value = 400.806900043077
$ws.Cells.Item(32, 4).Value = 400.806900043077
$ws.Cells.Item(32, 5).Value = 400.806900043077
$ws.Cells.Item(32, 6).Value = 400.806900043077
$ws.Cells.Item(32, 24).Value = 400.806900043077
$ws.Cells.Item(32, 25).Value = 400.806900043077
$ws.Cells.Item(35, 2).Value = 400.806900043077
$ws.Cells.Item(35, 3).Value = 400.806900043077
$ws.Cells.Item(35, 4).Value = 400.806900043077
$ws.Cells.Item(35, 5).Value = 400.806900043077
$ws.Cells.Item(35, 6).Value = 400.806900043077
$ws.Cells.Item(35, 24).Value = 400.806900043077
$ws.Cells.Item(35, 25).Value = 400.806900043077
$ws.Cells.Item(38, 2).Value = 400.806900043077
$ws.Cells.Item(38, 3).Value = 400.806900043077
$ws.Cells.Item(38, 4).Value = 400.806900043077
$ws.Cells.Item(38, 5).Value = 400.806900043077
$ws.Cells.Item(38, 6).Value = 400.806900043077
$ws.Cells.Item(38, 24).Value = 400.806900043077
$ws.Cells.Item(38, 25).Value = 400.806900043077
$ws.Cells.Item(41, 2).Value = 397.3838530629687
$ws.Cells.Item(41, 3).Value = 397.3838530629687
$ws.Cells.Item(41, 4).Value = 397.3838530629687
$ws.Cells.Item(41, 5).Value = 397.3838530629687
$ws.Cells.Item(41, 6).Value = 397.3838530629687
$ws.Cells.Item(41, 7).Value = 397.3838530629687
$ws.Cells.Item(41, 22).Value = 358.9907805655117
$ws.Cells.Item(41, 23).Value = 397.3838530629687
$ws.Cells.Item(41, 24).Value = 397.3838530629687
$ws.Cells.Item(41, 25).Value = 397.3838530629687
$ws.Cells.Item(44, 2).Value = 397.3838530629687
$ws.Cells.Item(44, 3).Value = 397.3838530629687
$ws.Cells.Item(44, 4).Value = 397.3838530629687
$ws.Cells.Item(44, 5).Value = 397.3838530629687
$ws.Cells.Item(44, 6).Value = 397.3838530629687
$ws.Cells.Item(44, 7).Value = 397.3838530629687
$ws.Cells.Item(44, 8).Value = 286.2388530112159
$ws.Cells.Item(44, 20).Value = 217.8665548556918
$ws.Cells.Item(44, 21).Value = 256.6300796561533
$ws.Cells.Item(44, 22).Value = 358.9907805655117
$ws.Cells.Item(44, 23).Value = 397.3838530629687
$ws.Cells.Item(44, 24).Value = 397.3838530629687
$ws.Cells.Item(44, 25).Value = 397.3838530629687
$ws.Cells.Item(45, 21).Value = 174.5731815300314
$ws.Cells.Item(45, 22).Value = 197.1263427586206
$ws.Cells.Item(45, 23).Value = 183.4695267241379
$ws.Cells.Item(45, 24).Value = 153.3187614035088
$ws.Cells.Item(46, 2).Value = 189.9004325317972
$ws.Cells.Item(46, 3).Value = 170.8360944016073
$ws.Cells.Item(46, 4).Value = 164.2192128704925
$ws.Cells.Item(46, 5).Value = 168.0604237117701
$ws.Cells.Item(46, 6).Value = 174.9399834978613
$ws.Cells.Item(46, 7).Value = 163.9353622244306
$ws.Cells.Item(46, 8).Value = 138.5031525665292
$ws.Cells.Item(46, 19).Value = 157.6489550149833
$ws.Cells.Item(46, 20).Value = 243.4206519573293
$ws.Cells.Item(46, 21).Value = 275.6486707394257
$ws.Cells.Item(46, 22).Value = 284.0859530482738
$ws.Cells.Item(46, 23).Value = 269.3061403695714
$ws.Cells.Item(46, 24).Value = 242.9378371199217
$ws.Cells.Item(46, 25).Value = 225.1454739790328

# --- Unmet Demand ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Cells.Item(11, 2).Value = 24.65272171708995
$ws.Cells.Item(11, 3).Value = 36.37719502544218
$ws.Cells.Item(11, 4).Value = 34.16683391433116
$ws.Cells.Item(11, 5).Value = 32.05314433031901
$ws.Cells.Item(11, 6).Value = 26.20490223191558
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 22).Value = 0
$ws.Cells.Item(11, 23).Value = 3.423046980108268
$ws.Cells.Item(11, 24).Value = 17.56718572351366
$ws.Cells.Item(11, 25).Value = 6.819409400374354
$ws.Cells.Item(17, 2).Value = 21.22967473698168
$ws.Cells.Item(17, 3).Value = 32.95414804533391
$ws.Cells.Item(17, 4).Value = 30.7437869342229
$ws.Cells.Item(17, 5).Value = 28.63009735021075
$ws.Cells.Item(17, 6).Value = 22.78185525180731
$ws.Cells.Item(17, 24).Value = 14.14413874340539
$ws.Cells.Item(17, 25).Value = 3.396362420266087
$ws.Cells.Item(20, 2).Value = 21.22967473698168
$ws.Cells.Item(20, 3).Value = 32.95414804533391
$ws.Cells.Item(20, 4).Value = 30.7437869342229
$ws.Cells.Item(20, 5).Value = 28.63009735021075
$ws.Cells.Item(20, 6).Value = 22.78185525180731
$ws.Cells.Item(20, 24).Value = 14.14413874340539
$ws.Cells.Item(20, 25).Value = 3.396362420266087
$ws.Cells.Item(23, 2).Value = 21.22967473698168
$ws.Cells.Item(23, 3).Value = 32.95414804533391
$ws.Cells.Item(23, 4).Value = 30.7437869342229
$ws.Cells.Item(23, 5).Value = 28.63009735021075
$ws.Cells.Item(23, 6).Value = 22.78185525180731
$ws.Cells.Item(23, 24).Value = 14.14413874340539
$ws.Cells.Item(23, 25).Value = 3.396362420266087
$ws.Cells.Item(26, 2).Value = 24.65272171708995
$ws.Cells.Item(26, 3).Value = 36.37719502544218
$ws.Cells.Item(26, 4).Value = 34.16683391433116
$ws.Cells.Item(26, 5).Value = 32.05314433031901
$ws.Cells.Item(26, 6).Value = 26.20490223191558
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 22).Value = 0
$ws.Cells.Item(26, 23).Value = 3.423046980108268
$ws.Cells.Item(26, 24).Value = 17.56718572351366
$ws.Cells.Item(26, 25).Value = 6.819409400374354
$ws.Cells.Item(29, 2).Value = 24.65272171708995
$ws.Cells.Item(29, 3).Value = 36.37719502544218
$ws.Cells.Item(29, 4).Value = 34.16683391433116
$ws.Cells.Item(29, 5).Value = 32.05314433031901
$ws.Cells.Item(29, 6).Value = 26.20490223191558
$ws.Cells.Item(29, 23).Value = 3.423046980108268
$ws.Cells.Item(29, 24).Value = 17.56718572351366
$ws.Cells.Item(29, 25).Value = 6.819409400374354
$ws.Cells.Item(32, 2).Value = 21.22967473698168
$ws.Cells.Item(32, 3).Value = 32.95414804533391
$ws.Cells.Item(32, 4).Value = 30.7437869342229
$ws.Cells.Item(32, 5).Value = 28.63009735021075
$ws.Cells.Item(32, 6).Value = 22.78185525180731
$ws.Cells.Item(32, 24).Value = 14.14413874340539
$ws.Cells.Item(32, 25).Value = 3.396362420266087
$ws.Cells.Item(35, 2).Value = 21.22967473698168
$ws.Cells.Item(35, 3).Value = 32.95414804533391
$ws.Cells.Item(35, 4).Value = 30.7437869342229
$ws.Cells.Item(35, 5).Value = 28.63009735021075
$ws.Cells.Item(35, 6).Value = 22.78185525180731
$ws.Cells.Item(35, 24).Value = 14.14413874340539
$ws.Cells.Item(35, 25).Value = 3.396362420266087
$ws.Cells.Item(38, 2).Value = 21.22967473698168
$ws.Cells.Item(38, 3).Value = 32.95414804533391
$ws.Cells.Item(38, 4).Value = 30.7437869342229
$ws.Cells.Item(38, 5).Value = 28.63009735021075
$ws.Cells.Item(38, 6).Value = 22.78185525180731
$ws.Cells.Item(38, 24).Value = 14.14413874340539
$ws.Cells.Item(38, 25).Value = 3.396362420266087
$ws.Cells.Item(41, 2).Value = 24.65272171708995
$ws.Cells.Item(41, 3).Value = 36.37719502544218
$ws.Cells.Item(41, 4).Value = 34.16683391433116
$ws.Cells.Item(41, 5).Value = 32.05314433031901
$ws.Cells.Item(41, 6).Value = 26.20490223191558
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 22).Value = 0
$ws.Cells.Item(41, 23).Value = 3.423046980108268
$ws.Cells.Item(41, 24).Value = 17.56718572351366
$ws.Cells.Item(41, 25).Value = 6.819409400374354
$ws.Cells.Item(44, 2).Value = 24.65272171708995
$ws.Cells.Item(44, 3).Value = 36.37719502544218
$ws.Cells.Item(44, 4).Value = 34.16683391433116
$ws.Cells.Item(44, 5).Value = 32.05314433031901
$ws.Cells.Item(44, 6).Value = 26.20490223191558
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 20).Value = 0
$ws.Cells.Item(44, 21).Value = 0
$ws.Cells.Item(44, 22).Value = 0
$ws.Cells.Item(44, 23).Value = 3.423046980108268
$ws.Cells.Item(44, 24).Value = 17.56718572351366
$ws.Cells.Item(44, 25).Value = 6.819409400374354
$ws.Cells.Item(45, 21).Value = 0
$ws.Cells.Item(45, 22).Value = 0
$ws.Cells.Item(45, 23).Value = 0
$ws.Cells.Item(45, 24).Value = 0
$ws.Cells.Item(46, 2).Value = 0
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 19).Value = 0
$ws.Cells.Item(46, 20).Value = 0
$ws.Cells.Item(46, 21).Value = 0
$ws.Cells.Item(46, 22).Value = 0
$ws.Cells.Item(46, 23).Value = 0
$ws.Cells.Item(46, 24).Value = 0
$ws.Cells.Item(46, 25).Value = 0

# --- Household Surplus ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Cells.Item(5, 2).Value = 391289.0028902176
$ws.Cells.Item(7, 2).Value = 392923.8501279173
$ws.Cells.Item(8, 2).Value = 392923.8501279173
$ws.Cells.Item(9, 2).Value = 392923.8501279173
$ws.Cells.Item(10, 2).Value = 391289.0028902176
$ws.Cells.Item(11, 2).Value = 391289.0028902176
$ws.Cells.Item(12, 2).Value = 392923.8501279173
$ws.Cells.Item(13, 2).Value = 392923.8501279173
$ws.Cells.Item(14, 2).Value = 392923.8501279173
$ws.Cells.Item(15, 2).Value = 391289.0028902176
$ws.Cells.Item(16, 2).Value = 391289.0028902176
